$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.502.10'
$ws.Range("E2").Value = '  -4.27%  '
$ws.Range("D3").Value = '2.974.20'
$ws.Range("E3").Value = '  -5.15%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.15'
$ws.Range("E5").Value = '  -5.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.68'
$ws.Range("E6").Value = '  -8.38%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -1.10%  '
$ws.Range("D9").Value = '2.981.60'
$ws.Range("E9").Value = '  -5.29%  '
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("E12").Value = '  -4.67%  '
$ws.Range("D13").Value = '3.489.00'
$ws.Range("E13").Value = '  -5.33%  '
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("D15").Value = '61.585.94'
$ws.Range("E15").Value = '  -4.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.61'
$ws.Range("E16").Value = '  -5.81%  '
$ws.Range("D17").Value = '2.975.14'
$ws.Range("E17").Value = '  -5.56%  '
$ws.Range("E18").Value = '  -5.23%  '
$ws.Range("E19").Value = '  -1.91%  '
$ws.Range("E20").Value = '  -3.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '379.30'
$ws.Range("E21").Value = '  -5.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.69'
$ws.Range("E22").Value = '  -5.30%  '
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.66'
$ws.Range("E24").Value = '  -3.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.67'
$ws.Range("E25").Value = '  -4.43%  '
$ws.Range("E26").Value = '  -2.94%  '
$ws.Range("D27").Value = '3.093.67'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.188'
$ws.Range("E28").Value = '  -4.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.996'
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = '0.0₃0937'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.21'
$ws.Range("E31").Value = '  -6.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.72'
$ws.Range("E33").Value = '  -5.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.41'
$ws.Range("E34").Value = '  -3.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.52'
$ws.Range("E35").Value = '  -1.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.64'
$ws.Range("E36").Value = '  -4.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.92'
$ws.Range("E37").Value = '  -5.53%  '
$ws.Range("E38").Value = '  -3.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.27'
$ws.Range("E39").Value = '  -5.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.55'
$ws.Range("E40").Value = '  -7.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.92'
$ws.Range("E41").Value = '  -3.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '37.56'
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("D43").Value = '2.410.74'
$ws.Range("E43").Value = '  -8.53%  '
$ws.Range("E44").Value = '  -6.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.670'
$ws.Range("E45").Value = '  -2.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0591'
$ws.Range("E46").Value = '  -3.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.996'
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.05'
$ws.Range("E48").Value = '  -6.88%  '
$ws.Range("E49").Value = '  -3.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0951'
$ws.Range("E50").Value = '  -2.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.73'
$ws.Range("E51").Value = '  -6.51%  '
